$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.153.13'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '2.544.37'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''591.59'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = '''173.76'
$ws.Range('E6').Value = '  +4.93%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = '2.545.45'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').Value = '''0.138'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('E11').Value = '  +1.99%  '
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '''0.346'
$ws.Range('E13').Value = '  -5.06%  '
$ws.Range('D14').Value = '''26.81'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('D15').Value = '3.009.77'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = '66.869.85'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '2.532.49'
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('D19').Value = '''8.09'
$ws.Range('E19').Value = '  +4.00%  '
$ws.Range('D20').Value = '''11.37'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').Value = '''355.03'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').Value = '''4.19'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').Value = '''4.62'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  +5.53%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '''69.91'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('D27').Value = '''10.11'
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '0.0₃0988'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = '''534.87'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').Value = '''8.17'
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '''1.46'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '''157.03'
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').Value = '''18.65'
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').Value = '''18.46'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').Value = '''0.356'
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('D42').Value = '''1.79'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '''5.14'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.51'
$ws.Range('E45').Value = '  +4.66%  '
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').Value = '''149.47'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').Value = '''0.561'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('D49').Value = '0.0₆0279'
$ws.Range('E49').Value = '  -4.20%  '
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('E51').Value = '  +0.01%  '
